$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.675229333333333
$ws.Range("H2").Value = 5.025688
$ws.Range("I2").Value = 0.2721044738138681
$ws.Range("J2").Value = 0.2721044738138681
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.301118333333333
$ws.Range("N2").Value = 3.903355
$ws.Range("O2").Value = 0.05808586029604726
$ws.Range("P2").Value = 0.05808586029604726
$ws.Range("Q2").Value = 2.179671598137778
$ws.Range("R2").Value = 19.61704438324
$ws.Range("S2").Value = 0.01580542245188179
$ws.Range("T2").Value = 0.01580542245188179
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.675229333333333
$ws.Range("H3").Value = 5.025688
$ws.Range("I3").Value = 0.2721044738138681
$ws.Range("J3").Value = 0.2721044738138681
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 11.20764866666667
$ws.Range("N3").Value = 33.622946
$ws.Range("O3").Value = 0.5003433569576788
$ws.Range("P3").Value = 0.500343356957679
$ws.Range("Q3").Value = 18.77538180409422
$ws.Range("R3").Value = 168.978436236848
$ws.Range("S3").Value = 0.1361456658712336
$ws.Range("T3").Value = 0.1361456658712336
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.675229333333333
$ws.Range("H4").Value = 5.025688
$ws.Range("I4").Value = 0.2721044738138681
$ws.Range("J4").Value = 0.2721044738138681
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.84722
$ws.Range("N4").Value = 11.54166
$ws.Range("O4").Value = 0.1717515445929148
$ws.Range("P4").Value = 0.1717515445929148
$ws.Range("Q4").Value = 6.444975795786666
$ws.Range("R4").Value = 58.00478216208
$ws.Range("S4").Value = 0.04673436366817417
$ws.Range("T4").Value = 0.04673436366817418
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.675229333333333
$ws.Range("H5").Value = 5.025688
$ws.Range("I5").Value = 0.2721044738138681
$ws.Range("J5").Value = 0.2721044738138681
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 6.043928
$ws.Range("N5").Value = 18.131784
$ws.Range("O5").Value = 0.2698192381533591
$ws.Range("P5").Value = 0.2698192381533591
$ws.Range("Q5").Value = 10.12496547415467
$ws.Range("R5").Value = 91.124689267392
$ws.Range("S5").Value = 0.07341902182257853
$ws.Range("T5").Value = 0.07341902182257853
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.433537333333333
$ws.Range("H6").Value = 7.300612
$ws.Range("I6").Value = 0.3952750721451891
$ws.Range("J6").Value = 0.3952750721451891
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.301118333333333
$ws.Range("N6").Value = 3.903355
$ws.Range("O6").Value = 0.05808586029604726
$ws.Range("P6").Value = 0.05808586029604726
$ws.Range("Q6").Value = 3.166320039251111
$ws.Range("R6").Value = 28.49688035326
$ws.Range("S6").Value = 0.02295989261913546
$ws.Range("T6").Value = 0.02295989261913546
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.433537333333333
$ws.Range("H7").Value = 7.300612
$ws.Range("I7").Value = 0.3952750721451891
$ws.Range("J7").Value = 0.3952750721451891
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 11.20764866666667
$ws.Range("N7").Value = 33.622946
$ws.Range("O7").Value = 0.5003433569576788
$ws.Range("P7").Value = 0.500343356957679
$ws.Range("Q7").Value = 27.27423144921688
$ws.Range("R7").Value = 245.468083042952
$ws.Range("S7").Value = 0.1977732565188126
$ws.Range("T7").Value = 0.1977732565188127
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.433537333333333
$ws.Range("H8").Value = 7.300612
$ws.Range("I8").Value = 0.3952750721451891
$ws.Range("J8").Value = 0.3952750721451891
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.84722
$ws.Range("N8").Value = 11.54166
$ws.Range("O8").Value = 0.1717515445929148
$ws.Range("P8").Value = 0.1717515445929148
$ws.Range("Q8").Value = 9.362353499546666
$ws.Range("R8").Value = 84.26118149592
$ws.Range("S8").Value = 0.06788910418001205
$ws.Range("T8").Value = 0.06788910418001207
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.433537333333333
$ws.Range("H9").Value = 7.300612
$ws.Range("I9").Value = 0.3952750721451891
$ws.Range("J9").Value = 0.3952750721451891
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 6.043928
$ws.Range("N9").Value = 18.131784
$ws.Range("O9").Value = 0.2698192381533591
$ws.Range("P9").Value = 0.2698192381533591
$ws.Range("Q9").Value = 14.70812442797867
$ws.Range("R9").Value = 132.373119851808
$ws.Range("S9").Value = 0.106652818827229
$ws.Range("T9").Value = 0.106652818827229
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.9268273333333333
$ws.Range("H10").Value = 2.780482
$ws.Range("I10").Value = 0.1505428891644152
$ws.Range("J10").Value = 0.1505428891644152
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.301118333333333
$ws.Range("N10").Value = 3.903355
$ws.Range("O10").Value = 0.05808586029604726
$ws.Range("P10").Value = 0.05808586029604726
$ws.Range("Q10").Value = 1.205912035234445
$ws.Range("R10").Value = 10.85320831711
$ws.Range("S10").Value = 0.00874441322856755
$ws.Range("T10").Value = 0.00874441322856755
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.9268273333333333
$ws.Range("H11").Value = 2.780482
$ws.Range("I11").Value = 0.1505428891644152
$ws.Range("J11").Value = 0.1505428891644152
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 11.20764866666667
$ws.Range("N11").Value = 33.622946
$ws.Range("O11").Value = 0.5003433569576788
$ws.Range("P11").Value = 0.500343356957679
$ws.Range("Q11").Value = 10.38755512666356
$ws.Range("R11").Value = 93.487996139972
$ws.Range("S11").Value = 0.0753231345306313
$ws.Range("T11").Value = 0.07532313453063132
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.9268273333333333
$ws.Range("H12").Value = 2.780482
$ws.Range("I12").Value = 0.1505428891644152
$ws.Range("J12").Value = 0.1505428891644152
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.84722
$ws.Range("N12").Value = 11.54166
$ws.Range("O12").Value = 0.1717515445929148
$ws.Range("P12").Value = 0.1717515445929148
$ws.Range("Q12").Value = 3.565708653346667
$ws.Range("R12").Value = 32.09137788012
$ws.Range("S12").Value = 0.02585597374146829
$ws.Range("T12").Value = 0.02585597374146829
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.9268273333333333
$ws.Range("H13").Value = 2.780482
$ws.Range("I13").Value = 0.1505428891644152
$ws.Range("J13").Value = 0.1505428891644152
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.043928
$ws.Range("N13").Value = 18.131784
$ws.Range("O13").Value = 0.2698192381533591
$ws.Range("P13").Value = 0.2698192381533591
$ws.Range("Q13").Value = 5.601677671098667
$ws.Range("R13").Value = 50.415099039888
$ws.Range("S13").Value = 0.04061936766374809
$ws.Range("T13").Value = 0.04061936766374809
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.120972666666667
$ws.Range("H14").Value = 3.362918
$ws.Range("I14").Value = 0.1820775648765275
$ws.Range("J14").Value = 0.1820775648765275
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.301118333333333
$ws.Range("N14").Value = 3.903355
$ws.Range("O14").Value = 0.05808586029604726
$ws.Range("P14").Value = 0.05808586029604726
$ws.Range("Q14").Value = 1.458518087765555
$ws.Range("R14").Value = 13.12666278989
$ws.Range("S14").Value = 0.01057613199646246
$ws.Range("T14").Value = 0.01057613199646246
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.120972666666667
$ws.Range("H15").Value = 3.362918
$ws.Range("I15").Value = 0.1820775648765275
$ws.Range("J15").Value = 0.1820775648765275
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 11.20764866666667
$ws.Range("N15").Value = 33.622946
$ws.Range("O15").Value = 0.5003433569576788
$ws.Range("P15").Value = 0.500343356957679
$ws.Range("Q15").Value = 12.56346781293644
$ws.Range("R15").Value = 113.071210316428
$ws.Range("S15").Value = 0.09110130003700131
$ws.Range("T15").Value = 0.09110130003700136
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.120972666666667
$ws.Range("H16").Value = 3.362918
$ws.Range("I16").Value = 0.1820775648765275
$ws.Range("J16").Value = 0.1820775648765275
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 3.84722
$ws.Range("N16").Value = 11.54166
$ws.Range("O16").Value = 0.1717515445929148
$ws.Range("P16").Value = 0.1717515445929148
$ws.Range("Q16").Value = 4.312628462653334
$ws.Range("R16").Value = 38.81365616388
$ws.Range("S16").Value = 0.03127210300326024
$ws.Range("T16").Value = 0.03127210300326025
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 1.120972666666667
$ws.Range("H17").Value = 3.362918
$ws.Range("I17").Value = 0.1820775648765275
$ws.Range("J17").Value = 0.1820775648765275
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 6.043928
$ws.Range("N17").Value = 18.131784
$ws.Range("O17").Value = 0.2698192381533591
$ws.Range("P17").Value = 0.2698192381533591
$ws.Range("Q17").Value = 6.775078087301333
$ws.Range("R17").Value = 60.975702785712
$ws.Range("S17").Value = 0.04912802983980345
$ws.Range("T17").Value = 0.04912802983980347
